$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.227022647857666
$ws.Range("B1").Value = 3.600754022598267
$ws.Range("C1").Value = 3.385593891143799
$ws.Range("D1").Value = 2.657873153686523
$ws.Range("E1").Value = 1.148515462875366
